$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Skill Description") - shifts old B (SFIA Level),
# C (Keycode), D (Description) one column to the right (C, D, E).
$ws.Columns.Item(2).Insert()

# Header for the new column
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Map each SkillCode (column A) to its full/friendly name for the new
# "Skill Description" column (B).
$map = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "ETDL"       = "Learning delivery"
    "MADE"       = "MADE"
    "INCA"       = "Content authoring"
    "ICPM"       = "Content publishing"
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -and $map.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $map[$code]
    }
}
